$wb = $excel.ActiveWorkbook

# "Repayment Schedule" sheet becomes the active sheet/tab (was "Transactions").
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

# Insert a new (blank spacer) column before column N, pushing the old
# N/O/P ("Late" / "Heading" / "Outstanding") columns one to the right.
$ws.Columns("N").Insert()

# Leave the selection where the author left it after the edit.
$ws.Range("S6").Select() | Out-Null
